$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.411.71'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '1.795.32'
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.58'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3799'
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("E8").Value = '  +1.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.81'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.201'
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07517'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.06'
$ws.Range("E13").Value = '  +7.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.475'
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("D15").Value = '1.794.31'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.079'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001105'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06666'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.95'
$ws.Range("E19").Value = '  +2.36%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.531'
$ws.Range("E21").Value = '  +4.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.40'
$ws.Range("E22").Value = '  +3.81%  '
$ws.Range("D23").Value = '27.419.65'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.53'
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.415'
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.565'
$ws.Range("E26").Value = '  +6.08%  '
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.48'
$ws.Range("E28").Value = '  +9.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.43'
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").Value = '1.997.83'
$ws.Range("E30").Value = '  +3.06%  '
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.048'
$ws.Range("E32").Value = '  -1.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.104'
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08702'
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.28'
$ws.Range("E35").Value = '  +2.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.654'
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6920'
$ws.Range("E37").Value = '  +8.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.458'
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06404'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.890'
$ws.Range("E40").Value = '  +3.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2205'
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02351'
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.271'
$ws.Range("E43").Value = '  +3.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.49'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6452'
$ws.Range("E45").Value = '  +3.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.874'
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.138'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.34'
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.68'
$ws.Range("E51").Value = '  +1.82%  '
